$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Total" label and a SUM formula totalling the BOM unit-cost column.
$ws.Range("A26").Value = "Total"
$ws.Range("B26").Formula = "=SUM(B2:B24)"
$ws.Range("B26").NumberFormat = "`"$`"#,##0_);[Red]\(`"$`"#,##0\)"

# Move the active selection like the source workbook ends up with.
$ws.Range("B27").Select()
